$wb = $excel.ActiveWorkbook
$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
# for the recalculated report date (08-Sep-2025 -> 16-Sep-2025, i.e. -8 days
# to expiry for every row).
# ---------------------------------------------------------------------------
$periodToExpire = @{
    3  = 426
    4  = 231
    5  = 239
    6  = 507
    7  = 358
    8  = 525
    9  = 247
    10 = 238
    11 = 400
    12 = 496
    13 = 350
    14 = 335
    15 = 504
    16 = 357
    17 = -22
    18 = -103
    19 = -126
    20 = -48
    21 = -48
    22 = 155
}

foreach ($row in $periodToExpire.Keys) {
    $wsTraining.Range("H$row").Value = $periodToExpire[$row]
    # Leading apostrophe keeps this a literal text label (not an auto-parsed
    # date serial) to match the report's "dd-MMM-yyyy" text column.
    $wsTraining.Range("I$row").Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# Exam Dashboard: comments column now reports date validity, and the column
# no longer needs to be as wide since the text is much shorter.
# ---------------------------------------------------------------------------
$wsExam.Columns.Item(5).ColumnWidth = 14.17

for ($row = 3; $row -le 10; $row++) {
    $wsExam.Range("E$row").Value = "date is valid"
}

# ---------------------------------------------------------------------------
# Styling: header rows (dark blue fill) switch from plain bold text to bold
# white text, and the dashboard titles drop their oversized 14pt font to
# match (now bold + white, default size).
# ---------------------------------------------------------------------------
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A1").Font.Color = 16777215
$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A1").Font.Color = 16777215

$wsTraining.Range("A2:K2").Font.Color = 16777215
$wsExam.Range("A2:G2").Font.Color = 16777215
